# Update column G ("K" - strikeouts) values per row to reflect
# regenerated save_data that uses K (strikeouts) instead of the old
# "Strike#" proxy metric. Only the G column values change; everything
# else on the sheet (headers, other stat columns, styles) stays intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 1;
    3 = 1;
    4 = 1;
    5 = 0;
    6 = 0;
    7 = 2;
    8 = 1;
    9 = 1;
    10 = 1;
    11 = 0;
    12 = 2;
    13 = 1;
    14 = 0;
    15 = 2;
    16 = 1;
    17 = 0;
    18 = 1;
    19 = 1;
    20 = 0;
    21 = 1;
    22 = 1;
    23 = 1;
    24 = 0;
    25 = 1;
    26 = 1;
    27 = 1;
    28 = 1;
    29 = 1;
    30 = 1;
    31 = 1;
    32 = 2;
    33 = 1;
    34 = 0;
    35 = 0;
    36 = 1;
    37 = 0;
    38 = 1;
    39 = 1;
    40 = 3;
    41 = 1;
    42 = 3;
    43 = 1;
    44 = 1;
    45 = 0;
    46 = 1;
    47 = 2;
    48 = 1;
    49 = 1;
    50 = 0;
    51 = 0;
    52 = 1;
    53 = 1;
    54 = 2;
    55 = 2;
    56 = 0;
    59 = 1;
    60 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
